# Auto-generated Excel COM-interop script applying the crypto price/volume refresh
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '42.662.68'
$ws.Range('E2').Value = '  -2.26%  '
$ws.Range('D3').Value = '2.348.46'
$ws.Range('E3').Value = '  -0.67%  '
$ws.Range('E4').Value = '  -0.26%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '325.03'
$ws.Range('E5').Value = '  +3.16%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '101.50'
$ws.Range('E6').Value = '  -6.38%  '
$ws.Range('E7').Value = '  -0.73%  '
$ws.Range('E8').Value = '  +0.09%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.621'
$ws.Range('E9').Value = '  -2.98%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '39.95'
$ws.Range('E10').Value = '  -7.23%  '
$ws.Range('E11').Value = '  -2.13%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '8.41'
$ws.Range('E12').Value = '  -4.22%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.994'
$ws.Range('E13').Value = '  -4.15%  '
$ws.Range('E14').Value = '  -0.05%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '16.08'
$ws.Range('E15').Value = '  -2.75%  '
$ws.Range('D16').Value = '2.705.70'
$ws.Range('E16').Value = '  -0.78%  '
$ws.Range('D17').Value = '2.349.37'
$ws.Range('E17').Value = '  -1.06%  '
$ws.Range('B18').Value = 'WrappedBTC'
$ws.Range('C18').Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range('D18').Value = '42.645.26'
$ws.Range('E18').Value = '  -2.26%  '
$ws.Range('B19').Value = 'Uniswap'
$ws.Range('C19').Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '7.91'
$ws.Range('E19').Value = '  +8.68%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '0.0000107'
$ws.Range('E20').Value = '  -2.38%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '76.41'
$ws.Range('E21').Value = '  +1.48%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '3.71'
$ws.Range('E22').Value = '  +7.22%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '263.78'
$ws.Range('E23').Value = '  +2.34%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '2.30'
$ws.Range('E24').Value = '  -9.39%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '10.06'
$ws.Range('E25').Value = '  +7.74%  '
$ws.Range('E26').Value = '  +0.12%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '11.43'
$ws.Range('E27').Value = '  -5.25%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '22.56'
$ws.Range('E28').Value = '  -0.81%  '
$ws.Range('E29').Value = '  -1.80%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '175.50'
$ws.Range('E30').Value = '  +1.25%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '3.10'
$ws.Range('E31').Value = '  -3.53%  '
$ws.Range('E32').Value = '  -3.12%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '35.24'
$ws.Range('E33').Value = '  -10.03%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '6.05'
$ws.Range('E34').Value = '  +0.84%  '
$ws.Range('E35').Value = '  -0.64%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '4.56'
$ws.Range('E36').Value = '  -8.22%  '
$ws.Range('E37').Value = '  +3.88%  '
$ws.Range('E38').Value = '  -4.97%  '
$ws.Range('E39').Value = '  -8.75%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '2.82'
$ws.Range('E40').Value = '  -0.11%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.237'
$ws.Range('E41').Value = '  +1.72%  '
$ws.Range('E42').Value = '  -0.82%  '
$ws.Range('E43').Value = '  -3.03%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '121.25'
$ws.Range('E44').Value = '  +8.54%  '
$ws.Range('B45').Value = 'FirstDigitalUSD'
$ws.Range('C45').Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '1.00'
$ws.Range('E45').Value = '  -0.18%  '
$ws.Range('B46').Value = 'BitcoinSV'
$ws.Range('C46').Value = 'https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '92.71'
$ws.Range('E46').Value = '  +24.49%  '
$ws.Range('E47').Value = '  -7.74%  '
$ws.Range('E48').Value = '  -2.24%  '
$ws.Range('E49').Value = '  -1.54%  '
$ws.Range('E50').Value = '  -3.97%  '
$ws.Range('E51').Value = '  -0.30%  '
